$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits (shared-string content) ---------------------------------

# Row 34: Soudan du Sud -> Soudan du Sud*  (now flagged like other
# resource-rich countries, so it also picks up the highlighted-row style)
$ws.Range("B34").Value = "Soudan du Sud*"

# Row 48: Cabo Verde -> Cabo Verde*  (style already matches the
# highlighted pattern, only the label changes)
$ws.Range("B48").Value = "Cabo Verde*"

# Row 57: Nigeria* -> Nigeria  (no longer flagged, loses the highlight)
$ws.Range("B57").Value = "Nigeria"

# Source footnote text
$ws.Range("A104").Value = "Source : World Population Prospects : the 2022 Revision - United Nations Department Economic and Social Affairs Population Division : World Population Prospects : the 2022 Revision - United Nations Department of Economic and Social Affairs Population Division, Global Hunger Index from Welthungerhilfe and Concern Worldwide (2021 edition)."

# --- Row highlight swap (B:J format only, values untouched) -------------
# Row 34 gains the "flagged country" fill; copy it from another row
# (20) that already carries that exact formatting.
$ws.Range("B20:J20").Copy()
$ws.Range("B34:J34").PasteSpecial(-4122)

# Row 57 loses the fill; copy plain formatting from row 40.
$ws.Range("B40:J40").Copy()
$ws.Range("B57:J57").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Recalculated aggregate values ---------------------------------------

$ws.Range("C69").Value = 62.213457142857202
$ws.Range("D69").Value = 64.267200000000003
$ws.Range("E69").Value = 60.1438285714286
$ws.Range("F69").Value = 4.6690857142857203
$ws.Range("G69").Value = 40.538457142857197
$ws.Range("H69").Value = 58.247842857142899
$ws.Range("I69").Value = 160.00498571428599
$ws.Range("J69").Value = 28.12

$ws.Range("C77").Value = 73.362863636363699
$ws.Range("D77").Value = 76.701327272727298
$ws.Range("E77").Value = 70.189490909090907
$ws.Range("F77").Value = 2.0467818181818198
$ws.Range("G77").Value = 13.384136363636401
$ws.Range("H77").Value = 16.1237909090909
$ws.Range("I77").Value = 53.9308727272727
$ws.Range("J77").Value = 7.7272727272727302

$ws.Range("C80").Value = 64.780749999999998
$ws.Range("D80").Value = 67.058840000000004
$ws.Range("E80").Value = 62.605229999999999
$ws.Range("F80").Value = 4.0200899999999997
$ws.Range("G80").Value = 39.23874
$ws.Range("H80").Value = 55.57246
$ws.Range("I80").Value = 134.48186000000001
$ws.Range("J80").Value = 21.828571428571401

$ws.Range("C82").Value = 63.217309090909097
$ws.Range("D82").Value = 65.542970454545497
$ws.Range("E82").Value = 60.9516909090909
$ws.Range("F82").Value = 4.1256636363636403
$ws.Range("G82").Value = 41.057459090909099
$ws.Range("H82").Value = 58.375184090909102
$ws.Range("I82").Value = 148.054086363636
$ws.Range("J82").Value = 24.494871794871798

$ws.Range("C84").Value = 60.992775000000002
$ws.Range("D84").Value = 62.993337500000003
$ws.Range("E84").Value = 58.999366666666702
$ws.Range("F84").Value = 4.8623791666666696
$ws.Range("G84").Value = 49.171525000000003
$ws.Range("H84").Value = 72.292095833333406
$ws.Range("I84").Value = 174.822483333333
$ws.Range("J84").Value = 29.976190476190499

$ws.Range("C86").Value = 64.738190909090903
$ws.Range("D86").Value = 67.1283590909091
$ws.Range("E86").Value = 62.449872727272698
$ws.Range("F86").Value = 3.7389000000000001
$ws.Range("G86").Value = 36.839490909090898
$ws.Range("H86").Value = 50.633527272727299
$ws.Range("I86").Value = 130.340454545455
$ws.Range("J86").Value = 19.675000000000001

$ws.Range("C87").Value = 70.363161290322594
$ws.Range("D87").Value = 73.163445161290298
$ws.Range("E87").Value = 67.688570967741896
$ws.Range("F87").Value = 2.6050225806451599
$ws.Range("G87").Value = 22.300493548387099
$ws.Range("H87").Value = 27.224848387096799
$ws.Range("I87").Value = 68.569480645161306
$ws.Range("J87").Value = 16.552

$ws.Range("C89").Value = 72.896337777777802
$ws.Range("D89").Value = 76.271726666666694
$ws.Range("E89").Value = 69.729984444444497
$ws.Range("F89").Value = 2.0043822222222198
$ws.Range("G89").Value = 11.970595555555599
$ws.Range("H89").Value = 14.2350688888889
$ws.Range("I89").Value = 48.7849066666667
$ws.Range("J89").Value = 7.5617647058823501

$ws.Range("C90").Value = 79.370791525423698
$ws.Range("D90").Value = 82.057852542372899
$ws.Range("E90").Value = 76.778615254237295
$ws.Range("F90").Value = 1.65551525423729
$ws.Range("G90").Value = 4.8009627118644103
$ws.Range("H90").Value = 5.7516694915254298
$ws.Range("I90").Value = 23.084274576271198
$ws.Range("J90").Value = 6.0428571428571498

$ws.Range("C94").Value = 71.436635483871001
$ws.Range("D94").Value = 74.354425806451601
$ws.Range("E94").Value = 68.771432258064607
$ws.Range("F94").Value = 2.3970709677419402
$ws.Range("G94").Value = 17.1926548387097
$ws.Range("H94").Value = 20.926122580645199
$ws.Range("I94").Value = 66.019241935483905

$ws.Range("C97").Value = 61.089160526315801
$ws.Range("D97").Value = 63.125728947368401
$ws.Range("E97").Value = 59.101689473684203
$ws.Range("F97").Value = 4.5708131578947402
$ws.Range("G97").Value = 48.728302631578998
$ws.Range("H97").Value = 70.263497368421099
$ws.Range("I97").Value = 171.511494736842
$ws.Range("J97").Value = 28.096875000000001

$ws.Range("C98").Value = 69.137342857142897
$ws.Range("D98").Value = 71.913014285714297
$ws.Range("E98").Value = 66.523342857142893
$ws.Range("F98").Value = 2.8566190476190498
$ws.Range("G98").Value = 27.778542857142899
$ws.Range("H98").Value = 34.2696857142857
$ws.Range("I98").Value = 83.8205047619048
$ws.Range("J98").Value = 21.766666666666701

# Row 99 (unlabelled trailing summary row) no longer carries values.
$ws.Range("C99:J99").ClearContents()

Write-Host "edit applied"
